$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 385

$ws.Range("C1").Value = "Image Path"

for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    if ($bVal -ne "No Image") {
        $aVal = $ws.Cells.Item($r, 1).Value2
        $ws.Cells.Item($r, 3).Value = "Data/images/$aVal.png"
    }
}

$ws.Range("A1:C385").Select() | Out-Null
